# The commit ("added visualization jupyter notebook") accompanied a small
# manual tidy-up of the "out" worksheet that was done in Excel before the
# notebook was added to the repo: the first column (the long date/run-id
# strings, e.g. "2020-11-03-12_56_57") was widened by hand so the values are
# no longer clipped, and the selection was left on C4.
#
# Everything else in the recorded XML diff (new xmlns/mc:Ignorable
# namespaces, xr:uid/revisionPtr, fileVersion/build numbers, workbookView
# geometry, dxf/cellXfs re-ordering, sortState xmlns, theme xmlns="" …) is
# purely an artifact of the file having been re-saved by a newer Excel
# build; it carries no semantic content and is not something a user
# action reproduces, so it is intentionally left alone here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A (was a narrow auto "best fit" width) to a fixed, explicit
# custom width so the full date/run-id strings are visible.
$ws.Columns("A").ColumnWidth = 24.7

# Leave the cursor/selection on C4 (matches the saved sheetView selection).
$ws.Range("C4").Select()
